$d = $word.ActiveDocument

$replacements = @(
    @("2023-12-22 Friday", "2023-12-23 Saturday"),
    @("73×53=3869", "58×52=3016"),
    @("48×82=3936", "99×87=8613"),
    @("81×90=7290", "77×81=6237"),
    @("45×58=2610", "98×55=5390"),
    @("38×19=722", "64×70=4480"),
    @("46×36=1656", "74×95=7030"),
    @("37×94=3478", "34×65=2210"),
    @("93×75=6975", "94×50=4700"),
    @("37×16=592", "45×16=720"),
    @("83×51=4233", "39×21=819"),
    @("50×92=4600", "30×45=1350"),
    @("28×94=2632", "87×58=5046"),
    @("39×41=1599", "94×50=4700"),
    @("20×45=900", "66×14=924"),
    @("19×31=589", "72×80=5760"),
    @("28×11=308", "43×57=2451"),
    @("78×15=1170", "43×40=1720"),
    @("48×81=3888", "74×95=7030"),
    @("46×19=874", "12×39=468"),
    @("77×58=4466", "15×96=1440"),
    @("95×71=6745", "85×33=2805"),
    @("24×98=2352", "11×90=990"),
    @("70×58=4060", "89×56=4984"),
    @("51×45=2295", "30×82=2460"),
    @("45×52=2340", "15×18=270")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
